$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

$ws.Range("D2").Value = "64.299.29"

$ws.Range("D3").Value = "3.496.82"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "588.36"
$ws.Range("E5").Value = "  +0.32%  "

Set-TextValue "D6" "133.83"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E8").Value = "  -0.42%  "

Set-TextValue "D9" "7.65"
$ws.Range("E9").Value = "  +6.41%  "

$ws.Range("E10").Value = "  -0.25%  "

$ws.Range("E11").Value = "  +0.17%  "

$ws.Range("D12").Value = "4.091.84"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("E14").Value = "  -0.60%  "

$ws.Range("D15").Value = "3.496.62"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").Value = "64.235.43"
$ws.Range("E16").Value = "  -0.12%  "

Set-TextValue "D17" "24.84"
$ws.Range("E17").Value = "  -3.65%  "

$ws.Range("E18").Value = "  +0.36%  "

$ws.Range("E19").Value = "  -0.63%  "

Set-TextValue "D20" "13.53"
$ws.Range("E20").Value = "  -1.58%  "

Set-TextValue "D21" "386.00"
$ws.Range("E21").Value = "  -0.24%  "

Set-TextValue "D22" "0.579"
$ws.Range("E22").Value = "  +2.17%  "

$ws.Range("D23").Value = "3.635.87"
$ws.Range("E23").Value = "  +0.05%  "

Set-TextValue "D24" "74.30"
$ws.Range("E24").Value = "  +0.26%  "

Set-TextValue "D25" "1.00"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("E26").Value = "  +0.19%  "

$ws.Range("E27").Value = "  +2.01%  "

Set-TextValue "D28" "0.999"
$ws.Range("E28").Value = "  -0.16%  "

Set-TextValue "D29" "7.29"
$ws.Range("E29").Value = "  -1.93%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D30" "1.50"
$ws.Range("E30").Value = "  +1.47%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D31" "2.24"
$ws.Range("E31").Value = "  +0.60%  "

Set-TextValue "D32" "8.14"
$ws.Range("E32").Value = "  -1.59%  "

$ws.Range("E33").Value = "  +4.15%  "

$ws.Range("D34").Value = "3.525.06"
$ws.Range("E34").Value = "  +0.24%  "

$ws.Range("E35").Value = "  -0.01%  "

Set-TextValue "D36" "23.27"
$ws.Range("E36").Value = "  -0.94%  "

Set-TextValue "D37" "5.40"
$ws.Range("E37").Value = "  +3.65%  "

Set-TextValue "D38" "6.89"
$ws.Range("E38").Value = "  +0.28%  "

$ws.Range("E39").Value = "  -1.06%  "

Set-TextValue "D40" "164.53"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("E43").Value = "  +0.00%  "

Set-TextValue "D44" "4.39"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("E45").Value = "  +1.06%  "

Set-TextValue "D46" "24.28"
$ws.Range("E46").Value = "  -5.95%  "

Set-TextValue "D47" "1.65"
$ws.Range("E47").Value = "  -0.72%  "

$ws.Range("D48").Value = "2.429.78"
$ws.Range("E48").Value = "  -2.08%  "

$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D49" "0.921"
$ws.Range("E49").Value = "  +2.57%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "6.79"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("E51").Value = "  -1.18%  "

